# Arreglado un error de formato en el excel
#
# The sheet held a single percolation run (q = 0.1 .. 1.0 in column A,
# "Componentes Conexos" count in column B). The committed numbers were
# wrong (they were off by several orders of magnitude) and only one run
# had been pasted in. This fixes the first run's values and appends the
# four other runs that belong in the same report, each block separated
# from the next by one blank row, exactly like the rest of the workbook
# already does between data blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected column-B values for the run already present in rows 2-11
# (column A, q = 0.1 .. 1.0, already holds the right values and is left
# untouched).
$run1 = @(37, 23, 13, 8, 4, 3, 2, 2, 1, 1)
for ($i = 0; $i -lt $run1.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $run1[$i]
}

# The four additional runs to append below, each preceded by a blank
# separator row.
$runs = @(
    ,@(34, 25, 20, 8, 5, 3, 2, 1, 1, 1)
    ,@(43, 35, 24, 17, 13, 6, 4, 2, 2, 1)
    ,@(36, 26, 16, 7, 5, 1, 1, 1, 1, 1)
    ,@(37, 19, 12, 6, 4, 1, 1, 1, 1, 1)
)

$row = 12
foreach ($run in $runs) {
    # Blank separator row between blocks.
    $ws.Cells.Item($row, 1).Value = "'"
    $ws.Cells.Item($row, 2).Value = "'"
    $row++

    for ($i = 0; $i -lt $run.Length; $i++) {
        $ws.Cells.Item($row, 1).Value = [math]::Round(0.1 * ($i + 1), 1)
        $ws.Cells.Item($row, 2).Value = $run[$i]
        $row++
    }
}
